$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36728
$ws.Range("C2").Value = 5799.800999999999
$ws.Range("D2").Value = 30928.199
$ws.Range("B3").Value = 35269
$ws.Range("C3").Value = 5649.136500000001
$ws.Range("D3").Value = 29619.8635
$ws.Range("B4").Value = 33348
$ws.Range("C4").Value = 5554.8325
$ws.Range("D4").Value = 27793.1675
$ws.Range("B5").Value = 31700
$ws.Range("C5").Value = 5481.98
$ws.Range("D5").Value = 26218.02
$ws.Range("B6").Value = 30762
$ws.Range("C6").Value = 5456.8815
$ws.Range("D6").Value = 25305.1185
$ws.Range("B7").Value = 30953
$ws.Range("C7").Value = 5512.395
$ws.Range("D7").Value = 25440.605
$ws.Range("B8").Value = 31330
$ws.Range("C8").Value = 5765.477227722772
$ws.Range("D8").Value = 25564.52277227723
$ws.Range("B9").Value = 32595
$ws.Range("C9").Value = 6519.4465
$ws.Range("D9").Value = 26075.5535
$ws.Range("B10").Value = 36208
$ws.Range("C10").Value = 7901.5895
$ws.Range("D10").Value = 28306.4105
$ws.Range("B11").Value = 39069
$ws.Range("C11").Value = 13052.0005
$ws.Range("D11").Value = 26016.9995
$ws.Range("B12").Value = 39884
$ws.Range("C12").Value = 14886.7985
$ws.Range("D12").Value = 24997.2015
$ws.Range("B13").Value = 39218
$ws.Range("C13").Value = 14810.005
$ws.Range("D13").Value = 24407.995
$ws.Range("B14").Value = 38683
$ws.Range("C14").Value = 14699.3385
$ws.Range("D14").Value = 23983.6615
$ws.Range("B15").Value = 40878
$ws.Range("C15").Value = 15357.636
$ws.Range("D15").Value = 25520.364
$ws.Range("B16").Value = 41508
$ws.Range("C16").Value = 15432.207
$ws.Range("D16").Value = 26075.793
$ws.Range("B17").Value = 40766
$ws.Range("C17").Value = 15338.9075
$ws.Range("D17").Value = 25427.0925
$ws.Range("B18").Value = 38994
$ws.Range("C18").Value = 15792.80217625723
$ws.Range("D18").Value = 23201.19782374277
$ws.Range("B19").Value = 39323
$ws.Range("C19").Value = 15557.50171551809
$ws.Range("D19").Value = 23765.49828448191
$ws.Range("B20").Value = 41039
$ws.Range("C20").Value = 15152.49877462994
$ws.Range("D20").Value = 25886.50122537006
$ws.Range("B21").Value = 15333
$ws.Range("C21").Value = 13812.90902852661
$ws.Range("D21").Value = 1520.09097147339
$ws.Range("C22").Value = 12007.84432898735
$ws.Range("D22").Value = 68523.995
$ws.Range("C23").Value = 9523.143
$ws.Range("D23").Value = 66852.495
$ws.Range("C24").Value = 6375.7855
$ws.Range("D24").Value = 33754.2325
$ws.Range("C25").Value = 5494.996500000001
$ws.Range("D25").Value = 0
